$d = $word.ActiveDocument

$replacements = @(
    @{old="94×81="; new="88×56="},
    @{old="46×40="; new="86×91="},
    @{old="30×90="; new="97×46="},
    @{old="66×99="; new="25×33="},
    @{old="31×35="; new="23×78="},
    @{old="91×68="; new="30×37="},
    @{old="68×90="; new="66×23="},
    @{old="26×82="; new="25×69="},
    @{old="61×78="; new="21×12="},
    @{old="92×83="; new="33×75="},
    @{old="26×64="; new="85×59="},
    @{old="12×27="; new="79×39="},
    @{old="70×34="; new="22×84="},
    @{old="94×18="; new="94×29="},
    @{old="56×50="; new="57×62="},
    @{old="88×48="; new="90×80="},
    @{old="87×71="; new="94×97="},
    @{old="88×37="; new="14×57="},
    @{old="96×44="; new="28×23="},
    @{old="55×64="; new="18×46="},
    @{old="60×57="; new="29×50="},
    @{old="55×28="; new="45×26="},
    @{old="31×66="; new="66×75="},
    @{old="48×84="; new="47×98="},
    @{old="96×49="; new="57×26="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
